$d = $word.ActiveDocument

# The document currently ends with two empty paragraphs right before the
# section properties. The target adds one more empty paragraph followed
# by three paragraphs containing "localStorage", "redux" and "firebase".

# 1) Create the four new paragraph breaks at the very end of the document.
for ($i = 1; $i -le 4; $i++) {
    $d.Paragraphs.Last.Range.InsertParagraphAfter()
}

# 2) Fill the last three of those new paragraphs with their text, leaving
#    the first newly-created paragraph empty.
$n = $d.Paragraphs.Count
$d.Paragraphs.Item($n - 2).Range.InsertAfter("localStorage")
$d.Paragraphs.Item($n - 1).Range.InsertAfter("redux")
$d.Paragraphs.Item($n).Range.InsertAfter("firebase")
